$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11.93525832574852
$ws.Range("D2").Value = 5.429571842611967
$ws.Range("E2").Value = 12.03446326042932
$ws.Range("F2").Value = 31.37004501887866
$ws.Range("G2").Value = 3.644148989456125
$ws.Range("L2").Value = 8.57919795471617
$ws.Range("M2").Value = 25.39126443584699
$ws.Range("N2").Value = 17.31393162051411
$ws.Range("O2").Value = 27.57016682120584
$ws.Range("C3").Value = 11.97759758143747
$ws.Range("D3").Value = 5.455844981540283
$ws.Range("E3").Value = 12.11641472385976
$ws.Range("F3").Value = 31.06515876136416
$ws.Range("G3").Value = 3.647910404639262
$ws.Range("L3").Value = 8.611272607793442
$ws.Range("M3").Value = 24.66761677828075
$ws.Range("N3").Value = 17.06337177794763
$ws.Range("O3").Value = 27.40352765126595
$ws.Range("C4").Value = 12.00700453916249
$ws.Range("D4").Value = 5.472732336385659
$ws.Range("E4").Value = 12.16918024285138
$ws.Range("F4").Value = 30.88758570865865
$ws.Range("G4").Value = 3.65034037958693
$ws.Range("L4").Value = 8.631936159813733
$ws.Range("M4").Value = 24.21377934149506
$ws.Range("N4").Value = 16.90975663314302
$ws.Range("O4").Value = 27.30983398726907
$ws.Range("C5").Value = 12.01984043253751
$ws.Range("D5").Value = 5.479804564211344
$ws.Range("E5").Value = 12.19129918281472
$ws.Range("F5").Value = 30.81771881880619
$ws.Range("G5").Value = 3.651361017739017
$ws.Range("L5").Value = 8.640601358720666
$ws.Range("M5").Value = 24.02669774776117
$ws.Range("N5").Value = 16.8472874239371
$ws.Range("O5").Value = 27.27384965003399
$ws.Range("C6").Value = 12.02202314454609
$ws.Range("D6").Value = 5.480990422582456
$ws.Range("E6").Value = 12.19500929184647
$ws.Range("F6").Value = 30.80627025889054
$ws.Range("G6").Value = 3.651532333616196
$ws.Range("L6").Value = 8.642055007843974
$ws.Range("M6").Value = 23.99551154256441
$ws.Range("N6").Value = 16.83692445246094
$ws.Range("O6").Value = 27.26800788818938
$ws.Range("C7").Value = 12.00717420457261
$ws.Range("D7").Value = 5.472826942979416
$ws.Range("E7").Value = 12.16947604854982
$ws.Range("F7").Value = 30.88663326116697
$ws.Range("G7").Value = 3.650354021015322
$ws.Range("L7").Value = 8.632052030139608
$ws.Range("M7").Value = 24.2112646073613
$ws.Range("N7").Value = 16.90891352817735
$ws.Range("O7").Value = 27.309339762065
$ws.Range("C8").Value = 11.94914591578256
$ws.Range("D8").Value = 5.438474313378568
$ws.Range("E8").Value = 12.06221300339991
$ws.Range("F8").Value = 31.26296833309964
$ws.Range("G8").Value = 3.645420995382662
$ws.Range("L8").Value = 8.590056555246859
$ws.Range("M8").Value = 25.14387208202385
$ws.Range("N8").Value = 17.2275341998564
$ws.Range("O8").Value = 27.51093797459647
$ws.Range("C9").Value = 11.86264450030498
$ws.Range("D9").Value = 5.377080441879198
$ws.Range("E9").Value = 11.87122847286348
$ws.Range("F9").Value = 32.07380959397845
$ws.Range("G9").Value = 3.63669781412018
$ws.Range("L9").Value = 8.515358879483529
$ws.Range("M9").Value = 26.88699463150682
$ws.Range("N9").Value = 17.85102946993977
$ws.Range("O9").Value = 27.97326460670584
$ws.Range("C10").Value = 11.81603881075373
$ws.Range("D10").Value = 5.335581790995248
$ws.Range("E10").Value = 11.74263132090609
$ws.Range("F10").Value = 32.70885882819673
$ws.Range("G10").Value = 3.630860970883794
$ws.Range("L10").Value = 8.465091868036911
$ws.Range("M10").Value = 28.10301845277555
$ws.Range("N10").Value = 18.30419320680317
$ws.Range("O10").Value = 28.35160154755031
$ws.Range("C11").Value = 11.79857987206242
$ws.Range("D11").Value = 5.317479398027938
$ws.Range("E11").Value = 11.68665687256163
$ws.Range("F11").Value = 33.00513969668998
$ws.Range("G11").Value = 3.628328290401964
$ws.Range("L11").Value = 8.443214402365001
$ws.Range("M11").Value = 28.63992687443359
$ws.Range("N11").Value = 18.50844512698667
$ws.Range("O11").Value = 28.53157816313827
$ws.Range("C12").Value = 11.79251169925263
$ws.Range("D12").Value = 5.310735525794627
$ws.Range("E12").Value = 11.66582274628739
$ws.Range("F12").Value = 33.11829619191692
$ws.Range("G12").Value = 3.627386728869665
$ws.Range("L12").Value = 8.435071362152982
$ws.Range("M12").Value = 28.84073671556652
$ws.Range("N12").Value = 18.58545086434456
$ws.Range("O12").Value = 28.60081285712857
$ws.Range("C13").Value = 11.79379435417784
$ws.Range("D13").Value = 5.312183004067961
$ws.Range("E13").Value = 11.67029365409148
$ws.Range("F13").Value = 33.0938847940621
$ws.Range("G13").Value = 3.627588733993269
$ws.Range("L13").Value = 8.436818831385832
$ws.Range("M13").Value = 28.79760263992916
$ws.Range("N13").Value = 18.56888247127682
$ws.Range("O13").Value = 28.58585465301761
$ws.Range("C14").Value = 11.79806972622347
$ws.Range("D14").Value = 5.316922351087369
$ws.Range("E14").Value = 11.68493558287015
$ws.Range("F14").Value = 33.01443047044371
$ws.Range("G14").Value = 3.628250477277859
$ws.Range("L14").Value = 8.442541638460028
$ws.Range("M14").Value = 28.65649873374004
$ws.Range("N14").Value = 18.51478762759976
$ws.Range("O14").Value = 28.5372527649446
$ws.Range("C15").Value = 11.80075939489653
$ws.Range("D15").Value = 5.319839795901884
$ws.Range("E15").Value = 11.6939513274381
$ws.Range("F15").Value = 32.96588455555987
$ws.Range("G15").Value = 3.628658091159744
$ws.Range("L15").Value = 8.446065425225168
$ws.Range("M15").Value = 28.56973754141081
$ws.Range("N15").Value = 18.48160673181994
$ws.Range("O15").Value = 28.5076220127243
$ws.Range("C16").Value = 11.81725552992135
$ws.Range("D16").Value = 5.336780395499911
$ws.Range("E16").Value = 11.74634011380615
$ws.Range("F16").Value = 32.68963628355682
$ws.Range("G16").Value = 3.631028943698217
$ws.Range("L16").Value = 8.466541452052647
$ws.Range("M16").Value = 28.06758826776447
$ws.Range("N16").Value = 18.29080090056558
$ws.Range("O16").Value = 28.33999377377613
$ws.Range("C17").Value = 11.82833747719443
$ws.Range("D17").Value = 5.347371244683264
$ws.Range("E17").Value = 11.77912499616146
$ws.Range("F17").Value = 32.52198864726846
$ws.Range("G17").Value = 3.632514689170423
$ws.Range("L17").Value = 8.479355651617823
$ws.Range("M17").Value = 27.75524569042458
$ws.Range("N17").Value = 18.17321352519779
$ws.Range("O17").Value = 28.2391399247215
$ws.Range("C18").Value = 11.8350634667643
$ws.Range("D18").Value = 5.353535843494738
$ws.Range("E18").Value = 11.79821971436703
$ws.Range("F18").Value = 32.42626431349947
$ws.Range("G18").Value = 3.633380790203238
$ws.Range("L18").Value = 8.486819199099733
$ws.Range("M18").Value = 27.5740722341807
$ws.Range("N18").Value = 18.10540593524556
$ws.Range("O18").Value = 28.18187532250173
$ws.Range("C19").Value = 11.83740106244978
$ws.Range("D19").Value = 5.355635625156292
$ws.Range("E19").Value = 11.80472571521484
$ws.Range("F19").Value = 32.39397737381562
$ws.Range("G19").Value = 3.633676022411338
$ws.Range("L19").Value = 8.489362251708336
$ws.Range("M19").Value = 27.51247391754177
$ws.Range("N19").Value = 18.08241959339237
$ws.Range("O19").Value = 28.16261575729846
$ws.Range("C20").Value = 11.82712132291223
$ws.Range("D20").Value = 5.346236276188089
$ws.Range("E20").Value = 11.77561039302095
$ws.Range("F20").Value = 32.53976308774367
$ws.Range("G20").Value = 3.632355335528467
$ws.Range("L20").Value = 8.477981922220907
$ws.Range("M20").Value = 27.78865390004353
$ws.Range("N20").Value = 18.18574945434336
$ws.Range("O20").Value = 28.24979935722273
$ws.Range("C21").Value = 11.7967991658206
$ws.Range("D21").Value = 5.315527277681861
$ws.Range("E21").Value = 11.68062507101183
$ws.Range("F21").Value = 33.03774285819275
$ws.Range("G21").Value = 3.628055632902114
$ws.Range("L21").Value = 8.440856877669443
$ws.Range("M21").Value = 28.69801362020035
$ws.Range("N21").Value = 18.53068634804645
$ws.Range("O21").Value = 28.55149937883699
$ws.Range("C22").Value = 11.78014996457794
$ws.Range("D22").Value = 5.29610463017602
$ws.Range("E22").Value = 11.62065735512091
$ws.Range("F22").Value = 33.3687570117561
$ws.Range("G22").Value = 3.625347537327166
$ws.Range("L22").Value = 8.417417820334638
$ws.Range("M22").Value = 29.27766188414861
$ws.Range("N22").Value = 18.75411119198308
$ws.Range("O22").Value = 28.75495643876038
$ws.Range("C23").Value = 11.78874446612252
$ws.Range("D23").Value = 5.306411753117604
$ws.Range("E23").Value = 11.65247041717529
$ws.Range("F23").Value = 33.19161412699031
$ws.Range("G23").Value = 3.626783601299659
$ws.Range("L23").Value = 8.429852516339716
$ws.Range("M23").Value = 28.96968581052393
$ws.Range("N23").Value = 18.63507080572995
$ws.Range("O23").Value = 28.64581025412185
$ws.Range("C24").Value = 11.82767004142807
$ws.Range("D24").Value = 5.346749159379769
$ws.Range("E24").Value = 11.77719857800932
$ws.Range("F24").Value = 32.53172520433451
$ws.Range("G24").Value = 3.632427342154125
$ws.Range("L24").Value = 8.478602684667278
$ws.Range("M24").Value = 27.77355503035079
$ws.Range("N24").Value = 18.1800825942781
$ws.Range("O24").Value = 28.24497798879281
$ws.Range("C25").Value = 11.88309118923197
$ws.Range("D25").Value = 5.393053267533014
$ws.Range("E25").Value = 11.92083105346116
$ws.Range("F25").Value = 31.84716575658939
$ws.Range("G25").Value = 3.638956673436301
$ws.Range("L25").Value = 8.534752571611428
$ws.Range("M25").Value = 26.42596154911183
$ws.Range("N25").Value = 18.56888247127682
$ws.Range("O25").Value = 28.58585465301761
